$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nino34")

# --- Updated forecast values for existing init months 2024-06, 2024-07, 2024-08 ---
# (rows 19, 20, 21 - lead columns B..U, i.e. columns 2..21)

$row19 = @(0.252,0.067,-0.178,-0.327,-0.484,-0.659,-0.785,-0.796,-0.709,-0.611,-0.555,-0.554,-0.587,-0.604,-0.578,-0.54,-0.541,-0.599,-0.669,-0.6860000000000001)
$row20 = @(0.192,0.07099999999999999,-0.118,-0.304,-0.473,-0.5669999999999999,-0.5570000000000001,-0.481,-0.411,-0.385,-0.407,-0.46,-0.5,-0.497,-0.476,-0.485,-0.538,-0.601,-0.614,-0.5590000000000001)
$row21 = @(-0.051,-0.137,-0.301,-0.458,-0.551,-0.542,-0.468,-0.409,-0.404,-0.463,-0.5649999999999999,-0.642,-0.648,-0.623,-0.639,-0.718,-0.8110000000000001,-0.834,-0.759,-0.632)
$row22 = @(-0.236,-0.324,-0.465,-0.532,-0.505,-0.431,-0.383,-0.389,-0.452,-0.552,-0.628,-0.636,-0.614,-0.631,-0.707,-0.793,-0.8110000000000001,-0.737,-0.615,-0.484)

for ($i = 0; $i -lt $row19.Length; $i++) {
    $ws.Cells.Item(19, 2 + $i).Value = $row19[$i]
}
for ($i = 0; $i -lt $row20.Length; $i++) {
    $ws.Cells.Item(20, 2 + $i).Value = $row20[$i]
}
for ($i = 0; $i -lt $row21.Length; $i++) {
    $ws.Cells.Item(21, 2 + $i).Value = $row21[$i]
}

# --- New row 22 for init month 2024-09 ---
# Copy formatting from the row above (row 21) so the new row matches the
# existing table styling (bold/centered/bordered label cell, numeric format).
$ws.Range("A21:U21").Copy() | Out-Null
$ws.Range("A22:U22").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Cells.Item(22, 1).Value = "2024-09"
for ($i = 0; $i -lt $row22.Length; $i++) {
    $ws.Cells.Item(22, 2 + $i).Value = $row22[$i]
}
